$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.147215366363525
$ws.Range("B1").Value = 1.052607536315918
$ws.Range("C1").Value = 0.8628063797950745
$ws.Range("D1").Value = 0.8764113187789917
$ws.Range("E1").Value = 0.9741549491882324
